$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# current_phase: 1 -> 2
$ws.Range("D26").Value = 2

# last_action_date: updated timestamp (stored as text, same as other date cells)
$ws.Range("E26").Value = "2026-02-19T04:57:04.605141+00:00"

# replies_count: 0 -> 1
$ws.Range("I26").Value = 1

# replied_message_ids: [] -> [10251]
$ws.Range("M26").Value = "[10251]"
